$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 978.8889
$ws.Range("I58").Value = 115.71429
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 347.14287
$ws.Range("L58").Value = 12000
$ws.Range("M58").Value = -197.14287
$ws.Range("N58").Value = -12300

$ws.Range("H64").Value = 2703.75
$ws.Range("I64").Value = 2419.1304
$ws.Range("J64").Value = 3431.111
$ws.Range("K64").Value = 2419.1304
$ws.Range("L64").Value = 3431.111
$ws.Range("M64").Value = -2171.1304
$ws.Range("N64").Value = -3927.111

$ws.Range("H67").Value = 2703.75
$ws.Range("I67").Value = 2419.1304
$ws.Range("J67").Value = 3431.111
$ws.Range("K67").Value = 2419.1304
$ws.Range("L67").Value = 3431.111
$ws.Range("M67").Value = -1561.1304
$ws.Range("N67").Value = -5147.111

$ws.Range("H70").Value = 10859
$ws.Range("I70").Value = 977
$ws.Range("J70").Value = 13682.429
$ws.Range("K70").Value = 2931
$ws.Range("L70").Value = 41047.287
$ws.Range("M70").Value = -2661
$ws.Range("N70").Value = -41587.287

$ws.Range("H73").Value = 10859
$ws.Range("I73").Value = 977
$ws.Range("J73").Value = 13682.429
$ws.Range("K73").Value = 2931
$ws.Range("L73").Value = 41047.287
$ws.Range("M73").Value = -1995
$ws.Range("N73").Value = -42919.287

$ws.Range("H76").Value = 20410818
$ws.Range("J76").Value = 2692.2856
$ws.Range("L76").Value = 2692.2856
$ws.Range("N76").Value = -3322.2856

$ws.Range("H79").Value = 20410818
$ws.Range("J79").Value = 2692.2856
$ws.Range("L79").Value = 2692.2856
$ws.Range("N79").Value = -4876.2856

$ws.Range("H82").Value = 3678.4
$ws.Range("I82").Value = 1826.2858
$ws.Range("J82").Value = 8000
$ws.Range("K82").Value = 5478.857400000001
$ws.Range("L82").Value = 24000
$ws.Range("M82").Value = -5072.857400000001
$ws.Range("N82").Value = -24812

$ws.Range("H85").Value = 3678.4
$ws.Range("I85").Value = 1826.2858
$ws.Range("J85").Value = 8000
$ws.Range("K85").Value = 5478.857400000001
$ws.Range("L85").Value = 24000
$ws.Range("M85").Value = -4074.857400000001
$ws.Range("N85").Value = -26808

$ws.Range("H88").Value = 6061918
$ws.Range("I88").Value = 738
$ws.Range("J88").Value = 9092508
$ws.Range("K88").Value = 738
$ws.Range("L88").Value = 9092508
$ws.Range("M88").Value = -332
$ws.Range("N88").Value = -9093320

$ws.Range("H91").Value = 6061918
$ws.Range("I91").Value = 738
$ws.Range("J91").Value = 9092508
$ws.Range("K91").Value = 738
$ws.Range("L91").Value = 9092508
$ws.Range("M91").Value = 666
$ws.Range("N91").Value = -9095316

$ws.Range("H137").Value = 13598568
$ws.Range("I137").Value = 1077.1111
$ws.Range("J137").Value = 32921318
$ws.Range("K137").Value = 3231.3333
$ws.Range("L137").Value = 98763954
$ws.Range("M137").Value = -681.3333000000002
$ws.Range("N137").Value = -98769054

$ws.Range("H138").Value = 2389.6526
$ws.Range("I138").Value = 1563.4348
$ws.Range("J138").Value = 3165.2856
$ws.Range("K138").Value = 4690.3044
$ws.Range("L138").Value = 9495.856800000001
$ws.Range("M138").Value = 449.6956
$ws.Range("N138").Value = -19775.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 23590930
$ws.Range("I74").Value = 22222884
$ws.Range("J74").Value = 26669032
$ws.Range("K74").Value = 22222884
$ws.Range("L74").Value = 26669032
$ws.Range("M74").Value = -22222010
$ws.Range("N74").Value = -26670780

$ws.Range("H77").Value = 23590930
$ws.Range("I77").Value = 22222884
$ws.Range("J77").Value = 26669032
$ws.Range("K77").Value = 111114420
$ws.Range("L77").Value = 133345160
$ws.Range("M77").Value = -111110052
$ws.Range("N77").Value = -133353896

$ws.Range("H132").Value = 7938808
$ws.Range("I132").Value = 7938957.5
$ws.Range("J132").Value = 7937466.5
$ws.Range("K132").Value = 23816872.5
$ws.Range("L132").Value = 23812399.5
$ws.Range("M132").Value = -23814342.5
$ws.Range("N132").Value = -23817459.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 455.82758
$ws.Range("I22").Value = 467.37036
$ws.Range("K22").Value = 467.37036
$ws.Range("M22").Value = -294.37036

$ws.Range("H134").Value = 10505284
$ws.Range("I134").Value = 13514405
$ws.Range("J134").Value = 2552607
$ws.Range("K134").Value = 40543215
$ws.Range("L134").Value = 7657821
$ws.Range("M134").Value = -40540680
$ws.Range("N134").Value = -7662891

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 824830.9399999999
$ws.Range("I31").Value = 1127.8
$ws.Range("J31").Value = 3913717.8
$ws.Range("K31").Value = 1127.8
$ws.Range("L31").Value = 3913717.8
$ws.Range("M31").Value = -832.8
$ws.Range("N31").Value = -3914307.8

$ws.Range("H34").Value = 824830.9399999999
$ws.Range("I34").Value = 1127.8
$ws.Range("J34").Value = 3913717.8
$ws.Range("K34").Value = 1127.8
$ws.Range("L34").Value = 3913717.8
$ws.Range("M34").Value = -925.8
$ws.Range("N34").Value = -3914121.8

$ws.Range("H58").Value = 772881.5600000001
$ws.Range("I58").Value = 3258.359
$ws.Range("J58").Value = 2273646.8
$ws.Range("K58").Value = 3258.359
$ws.Range("L58").Value = 2273646.8
$ws.Range("M58").Value = -3055.359
$ws.Range("N58").Value = -2274052.8

$ws.Range("H134").Value = 600218.3
$ws.Range("I134").Value = 3203.673
$ws.Range("J134").Value = 2669869
$ws.Range("K134").Value = 9611.019
$ws.Range("L134").Value = 8009607
$ws.Range("M134").Value = -7076.019
$ws.Range("N134").Value = -8014677

$ws.Range("H136").Value = 772881.5600000001
$ws.Range("I136").Value = 3258.359
$ws.Range("J136").Value = 2273646.8
$ws.Range("K136").Value = 9775.076999999999
$ws.Range("L136").Value = 6820940.399999999
$ws.Range("M136").Value = -7225.076999999999
$ws.Range("N136").Value = -6826040.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 166751680
$ws.Range("I62").Value = 2888
$ws.Range("J62").Value = 333500480
$ws.Range("K62").Value = 8664
$ws.Range("L62").Value = 1000501440
$ws.Range("M62").Value = -7978
$ws.Range("N62").Value = -1000502812

$ws.Range("H65").Value = 166751680
$ws.Range("I65").Value = 2888
$ws.Range("J65").Value = 333500480
$ws.Range("K65").Value = 25992
$ws.Range("L65").Value = 3001504320
$ws.Range("M65").Value = -22560
$ws.Range("N65").Value = -3001511184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 15388371
$ws.Range("I132").Value = 12381693
$ws.Range("J132").Value = 45455144
$ws.Range("K132").Value = 37145079
$ws.Range("L132").Value = 136365432
$ws.Range("M132").Value = -37142549
$ws.Range("N132").Value = -136370492

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4087301.8
$ws.Range("I132").Value = 5296117.5
$ws.Range("J132").Value = 7548.125
$ws.Range("K132").Value = 15888352.5
$ws.Range("L132").Value = 22644.375
$ws.Range("M132").Value = -15885822.5
$ws.Range("N132").Value = -27704.375

$ws.Range("H136").Value = 1852729.1
$ws.Range("I136").Value = 1984999.1
$ws.Range("J136").Value = 950
$ws.Range("K136").Value = 5954997.300000001
$ws.Range("L136").Value = 2850
$ws.Range("M136").Value = -5952447.300000001
$ws.Range("N136").Value = -7950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 823260.75
$ws.Range("I132").Value = 2253.3044
$ws.Range("J132").Value = 3970455.8
$ws.Range("K132").Value = 6759.9132
$ws.Range("L132").Value = 11911367.4
$ws.Range("M132").Value = -4229.9132
$ws.Range("N132").Value = -11916427.4
